$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.367.75"
$ws.Range("E2").Value = "'  -4.86%  "
$ws.Range("D3").Value = "'1.569.55"
$ws.Range("E3").Value = "'  -4.90%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E5").Value = "'  -0.03%  "
$ws.Range("D6").Value = "'289.82"
$ws.Range("E6").Value = "'  -3.46%  "
$ws.Range("D7").Value = "'0.3770"
$ws.Range("E7").Value = "'  -0.34%  "
$ws.Range("D8").Value = "'49.64"
$ws.Range("E8").Value = "'  -2.51%  "
$ws.Range("E9").Value = "'  -2.57%  "
$ws.Range("D10").Value = "'1.166"
$ws.Range("E10").Value = "'  -4.93%  "
$ws.Range("D11").Value = "'0.07688"
$ws.Range("E11").Value = "'  -4.64%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "'  +0.06%  "
$ws.Range("D13").Value = "'21.39"
$ws.Range("E13").Value = "'  -2.88%  "
$ws.Range("D14").Value = "'6.035"
$ws.Range("E14").Value = "'  -4.69%  "
$ws.Range("E15").Value = "'  -4.47%  "
$ws.Range("D16").Value = "'0.00001139"
$ws.Range("E16").Value = "'  -5.83%  "
$ws.Range("D17").Value = "'1.571.34"
$ws.Range("E17").Value = "'  -4.70%  "
$ws.Range("D18").Value = "'90.59"
$ws.Range("E18").Value = "'  -4.86%  "
$ws.Range("E19").Value = "'  -3.89%  "
$ws.Range("E20").Value = "'  -0.05%  "
$ws.Range("D21").Value = "'6.258"
$ws.Range("E21").Value = "'  -5.74%  "
$ws.Range("D22").Value = "'16.67"
$ws.Range("D23").Value = "'0.5334"
$ws.Range("E23").Value = "'  -7.76%  "
$ws.Range("D24").Value = "'11.98"
$ws.Range("E24").Value = "'  -3.93%  "
$ws.Range("D25").Value = "'22.383.70"
$ws.Range("E25").Value = "'  -4.76%  "
$ws.Range("D26").Value = "'2.390"
$ws.Range("E26").Value = "'  -1.39%  "
$ws.Range("D27").Value = "'2.789"
$ws.Range("E27").Value = "'  -8.04%  "
$ws.Range("D28").Value = "'20.21"
$ws.Range("E28").Value = "'  -4.28%  "
$ws.Range("D29").Value = "'145.02"
$ws.Range("E29").Value = "'  -4.14%  "
$ws.Range("D30").Value = "'5.016"
$ws.Range("E30").Value = "'  -2.94%  "
$ws.Range("D31").Value = "'126.02"
$ws.Range("E31").Value = "'  -4.39%  "
$ws.Range("D32").Value = "'1.754.80"
$ws.Range("E32").Value = "'  -4.19%  "
$ws.Range("D33").Value = "'1.022"
$ws.Range("E33").Value = "'  +2.84%  "
$ws.Range("D34").Value = "'6.230"
$ws.Range("E34").Value = "'  -9.64%  "
$ws.Range("D35").Value = "'2.014"
$ws.Range("E35").Value = "'  -6.12%  "
$ws.Range("E36").Value = "'  -8.55%  "
$ws.Range("D37").Value = "'0.08533"
$ws.Range("E37").Value = "'  -2.95%  "
$ws.Range("D38").Value = "'0.02545"
$ws.Range("E38").Value = "'  -6.64%  "
$ws.Range("E39").Value = "'  -4.02%  "
$ws.Range("D40").Value = "'5.562"
$ws.Range("E40").Value = "'  -6.07%  "
$ws.Range("D41").Value = "'1.329"
$ws.Range("E41").Value = "'  +2.49%  "
$ws.Range("D42").Value = "'0.06398"
$ws.Range("E42").Value = "'  -7.06%  "
$ws.Range("D43").Value = "'11.72"
$ws.Range("E43").Value = "'  -8.93%  "
$ws.Range("D44").Value = "'0.6416"
$ws.Range("E44").Value = "'  -7.07%  "
$ws.Range("D45").Value = "'14.21"
$ws.Range("E45").Value = "'  -9.34%  "
$ws.Range("E46").Value = "'  -0.04%  "
$ws.Range("D47").Value = "'0.5999"
$ws.Range("E47").Value = "'  -6.14%  "
$ws.Range("D48").Value = "'3.759"
$ws.Range("E48").Value = "'  -4.12%  "
$ws.Range("E49").Value = "'  -6.95%  "
$ws.Range("D50").Value = "'1.306"
$ws.Range("E50").Value = "'  +5.22%  "
$ws.Range("D51").Value = "'124.54"
$ws.Range("E51").Value = "'  -2.02%  "
